$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Part 1: "Your report ... should be submitted as a PDF with a ..."
#         -> split the run and replace "PDF" with "hard copy to ITO"
# -----------------------------------------------------------------
$find = $d.Content.Find
$find.Execute("Your report for both parts should be submitted as a PDF with a cover page including just your name, student number and course details. Late submissions will not be accepted.", `
              $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$fullRange = $find.Parent
$start = $fullRange.Start
$end = $fullRange.End

$prefix = "Your report for both parts should be submitted as a "
$prefixLen = $prefix.Length
$oldMiddle = "PDF"
$oldMiddleLen = $oldMiddle.Length
$newMiddle = "hard copy to ITO"
$newMiddleLen = $newMiddle.Length

$midStart = $start + $prefixLen
$midEnd = $midStart + $oldMiddleLen

# Replace "PDF" with "hard copy to ITO"
$midRange = $d.Range($midStart, $midEnd)
$midRange.Text = $newMiddle

$newEnd = $end + ($newMiddleLen - $oldMiddleLen)

# Touch (no-op) formatting on each of the three logical spans so the
# engine keeps them as three distinct runs instead of re-merging them.
$r1 = $d.Range($start, $midStart)
$r1.Font.Bold = 1
$r1.Font.Bold = 0

$r2 = $d.Range($midStart, $midStart + $newMiddleLen)
$r2.Font.Bold = 1
$r2.Font.Bold = 0

$r3 = $d.Range($midStart + $newMiddleLen, $newEnd)
$r3.Font.Bold = 1
$r3.Font.Bold = 0

# -----------------------------------------------------------------
# Part 2: add the 19 new "ListLabel 115".."ListLabel 133" character
#         styles that appear in the saved styles.xml.
# -----------------------------------------------------------------
function Add-ListLabelStyle($id, $fontCs, $underline, $color) {
    $name = "ListLabel" + $id
    $s = $d.Styles.Add($name, 2)
    $s.NameLocal = "ListLabel " + $id
    $s.QuickStyle = $true
    if ($fontCs -ne $null) {
        $s.Font.NameBi = $fontCs
    }
    if ($color -ne $null) {
        $s.Font.Color = $color
    }
    $s.Font.Underline = $underline
}

# ListLabel115 .. ListLabel123 : plain, underline "none"
for ($i = 115; $i -le 123; $i++) {
    Add-ListLabelStyle $i $null 0 $null
}

# ListLabel124 .. ListLabel132 : complex-script font cycling through the
# three bullet fonts used by the document's bullet list levels, underline "none"
$fontsCycle = @("Wingdings", "Wingdings 2", "OpenSymbol")
$fi = 0
for ($i = 124; $i -le 132; $i++) {
    $font = $fontsCycle[$fi % 3]
    $fi = $fi + 1
    Add-ListLabelStyle $i $font 0 $null
}

# ListLabel133 : same look as ListLabel114 (hyperlink-like blue, single underline)
Add-ListLabelStyle 133 $null 1 13391121
